# Refresh the cryptocurrency price/volume table (rows 2-51) to the latest
# scraped values, matching the upstream "Updated cryptos list" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "29.334.93"
$ws.Range("E2").Formula = "  -0.46%  "

# Row 3
$ws.Range("D3").Formula = "1.845.65"
$ws.Range("E3").Formula = "  -0.34%  "

# Row 4
$ws.Range("D4").Formula = "'0.9981"
$ws.Range("E4").Formula = "  -0.18%  "

# Row 5
$ws.Range("D5").Formula = "'240.47"
$ws.Range("E5").Formula = "  -0.18%  "

# Row 6
$ws.Range("D6").Formula = "'0.6277"
$ws.Range("E6").Formula = "  -0.52%  "

# Row 7
$ws.Range("D7").Formula = "'0.9999"
$ws.Range("E7").Formula = "  -0.09%  "

# Row 8
$ws.Range("D8").Formula = "'0.07569"
$ws.Range("E8").Formula = "  -1.54%  "

# Row 9
$ws.Range("D9").Formula = "'0.2908"
$ws.Range("E9").Formula = "  -0.88%  "

# Row 10
$ws.Range("D10").Formula = "'24.46"
$ws.Range("E10").Formula = "  -1.23%  "

# Row 11
$ws.Range("D11").Formula = "'0.07751"
$ws.Range("E11").Formula = "  +0.01%  "

# Row 12
$ws.Range("D12").Formula = "1.846.37"
$ws.Range("E12").Formula = "  -1.12%  "

# Row 13
$ws.Range("D13").Formula = "'5.003"
$ws.Range("E13").Formula = "  -0.69%  "

# Row 14
$ws.Range("D14").Formula = "'0.6776"
$ws.Range("E14").Formula = "  -0.50%  "

# Row 15
$ws.Range("D15").Formula = "'0.00001037"
$ws.Range("E15").Formula = "  -3.20%  "

# Row 16
$ws.Range("D16").Formula = "'83.03"
$ws.Range("E16").Formula = "  -0.87%  "

# Row 17
$ws.Range("D17").Formula = "'6.097"
$ws.Range("E17").Formula = "  -1.70%  "

# Row 18
$ws.Range("D18").Formula = "29.321.81"
$ws.Range("E18").Formula = "  -0.56%  "

# Row 19
$ws.Range("D19").Formula = "'228.92"
$ws.Range("E19").Formula = "  -0.06%  "

# Row 20
$ws.Range("D20").Formula = "'12.31"
$ws.Range("E20").Formula = "  -1.27%  "

# Row 21
$ws.Range("D21").Formula = "'0.9994"
$ws.Range("E21").Formula = "  -0.09%  "

# Row 22
$ws.Range("D22").Formula = "'7.427"
$ws.Range("E22").Formula = "  -0.54%  "

# Row 23
$ws.Range("D23").Formula = "'1.002"
$ws.Range("E23").Formula = "  +0.15%  "

# Row 24
$ws.Range("D24").Formula = "'158.85"
$ws.Range("E24").Formula = "  +0.99%  "

# Row 25
$ws.Range("E25").Formula = "  +0.70%  "

# Row 26
$ws.Range("D26").Formula = "'8.430"
$ws.Range("E26").Formula = "  +0.17%  "

# Row 27
$ws.Range("D27").Formula = "'17.65"
$ws.Range("E27").Formula = "  -0.32%  "

# Row 28
$ws.Range("D28").Formula = "'1.415"
$ws.Range("E28").Formula = "  +5.94%  "

# Row 29
$ws.Range("D29").Formula = "'1.472"
$ws.Range("E29").Formula = "  +0.36%  "

# Row 30
$ws.Range("D30").Formula = "'0.05685"
$ws.Range("E30").Formula = "  -0.19%  "

# Row 31
$ws.Range("D31").Formula = "'4.104"
$ws.Range("E31").Formula = "  -0.74%  "

# Row 32
$ws.Range("D32").Formula = "'4.038"
$ws.Range("E32").Formula = "  -0.12%  "

# Row 33
$ws.Range("B33").Formula = "ARBITRUM"
$ws.Range("C33").Formula = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Formula = "'1.153"
$ws.Range("E33").Formula = "  -1.15%  "

# Row 34
$ws.Range("B34").Formula = "LidoDAOToken"
$ws.Range("C34").Formula = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Formula = "'1.818"
$ws.Range("E34").Formula = "  -1.94%  "

# Row 35
$ws.Range("D35").Formula = "'0.6963"
$ws.Range("E35").Formula = "  -1.75%  "

# Row 36
$ws.Range("D36").Formula = "'2.579"
$ws.Range("E36").Formula = "  -0.33%  "

# Row 37
$ws.Range("D37").Formula = "'0.01830"
$ws.Range("E37").Formula = "  +2.04%  "

# Row 38
$ws.Range("D38").Formula = "1.236.94"
$ws.Range("E38").Formula = "  +1.31%  "

# Row 39
$ws.Range("D39").Formula = "'2.717"
$ws.Range("E39").Formula = "  -2.41%  "

# Row 40
$ws.Range("D40").Formula = "'6.412"
$ws.Range("E40").Formula = "  -2.04%  "

# Row 41
$ws.Range("D41").Formula = "'0.9000"
$ws.Range("E41").Formula = "  -0.90%  "

# Row 42
$ws.Range("D42").Formula = "'0.9994"
$ws.Range("E42").Formula = "  -0.20%  "

# Row 43
$ws.Range("D43").Formula = "2.004.40"
$ws.Range("E43").Formula = "  -1.52%  "

# Row 44
$ws.Range("D44").Formula = "'101.43"
$ws.Range("E44").Formula = "  -0.54%  "

# Row 45
$ws.Range("D45").Formula = "'65.45"
$ws.Range("E45").Formula = "  -1.75%  "

# Row 46
$ws.Range("D46").Formula = "'7.122"
$ws.Range("E46").Formula = "  -0.15%  "

# Row 47
$ws.Range("B47").Formula = "TheSandbox"
$ws.Range("C47").Formula = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").Formula = "'0.3992"
$ws.Range("E47").Formula = "  -0.83%  "

# Row 48
$ws.Range("B48").Formula = "EnergySwap"
$ws.Range("C48").Formula = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Formula = "'9.005"
$ws.Range("E48").Formula = "  -0.46%  "

# Row 49
$ws.Range("B49").Formula = "Algorand"
$ws.Range("C49").Formula = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Formula = "'0.1148"
$ws.Range("E49").Formula = "  +0.34%  "

# Row 50
$ws.Range("B50").Formula = "BabyDogeCoin"
$ws.Range("C50").Formula = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Formula = "'0.00000000115"
$ws.Range("E50").Formula = "  -4.42%  "

# Row 51
$ws.Range("D51").Formula = "'1.671"
$ws.Range("E51").Formula = "  -1.07%  "
